$wb = $excel.ActiveWorkbook

# --- Leave the previously-active sheet with a whole-column selection, as   ---
# --- happened before the new sheet was created and made active.           ---
$ii = $wb.Worksheets.Item("indications_interventions")
$ii.Activate()
$ii.Columns("A:C").Select()

# --- Add the new "studyDesignPopulation" sheet after the last sheet ---
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "studyDesignPopulation"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 20.166666666666668
$ws.Columns.Item(2).ColumnWidth = 26.333333333333336
$ws.Columns.Item(3).ColumnWidth = 30.0
$ws.Columns.Item(4).ColumnWidth = 30.5
$ws.Columns.Item(5).ColumnWidth = 25.833333333333336

# --- Content (entered in the order that reproduces the shared-string table) ---
$ws.Range("A1").Value = "populationDescription"
$ws.Range("A2").Value = "Pop 1"
$ws.Range("B1").Value = "plannedNumberOfParticipants"
$ws.Range("B2").Value = 100
$ws.Range("C1").Value = "plannedMinimumAgeOfParticipants"
$ws.Range("C2").Value = "18 years"
$ws.Range("D2").Value = "'40 years"
$ws.Range("E1").Value = "plannedSexOfParticipants"
$ws.Range("E2").Value = "BOTH"

$ws.Range("A3").Value = "Pop 2"
$ws.Range("B3").Value = 20
$ws.Range("C3").Value = "18 years"
$ws.Range("D3").Value = "'60 years"
$ws.Range("E3").Value = "M"

$ws.Range("A4").Value = "Pop 3"
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = "18 years"
$ws.Range("D4").Value = "'70 years"
$ws.Range("E4").Value = "F"

$ws.Range("D1").Value = "plannedMaximumAgeOfParticipants"

# --- Bold header row (covers A1:G1, including the two trailing empty cells) ---
$ws.Range("A1:G1").Font.Bold = $true

# --- Make the new sheet the active tab, zoomed to 170%, with B8 selected ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 170
$ws.Range("B8").Select()
